$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12245
$ws1.Range("F4").Value = 4497
$ws1.Range("F5").Value = 49
$ws1.Range("F6").Value = 68
$ws1.Range("F8").Value = 30
$ws1.Range("F9").Value = 2613
$ws1.Range("F10").Value = 1134
$ws1.Range("F11").Value = 208
$ws1.Range("F12").Value = 79
$ws1.Range("F13").Value = 5342
$ws1.Range("F15").Value = 214
$ws1.Range("F16").Value = 562
$ws1.Range("F17").Value = 11494
$ws1.Range("F18").Value = 11602
$ws1.Range("F19").Value = 23
$ws1.Range("F20").Value = 66
$ws1.Range("F24").Value = 26

# Sheet "全部类型" (all types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 12245
$ws4.Range("F4").Value = 4497
$ws4.Range("F5").Value = 49
$ws4.Range("F6").Value = 68
$ws4.Range("F8").Value = 30
$ws4.Range("F9").Value = 2613
$ws4.Range("F11").Value = 1134
$ws4.Range("F12").Value = 208
$ws4.Range("F13").Value = 79
$ws4.Range("F14").Value = 5342
$ws4.Range("F16").Value = 214
$ws4.Range("F17").Value = 562
$ws4.Range("F18").Value = 11494
$ws4.Range("F19").Value = 11602
$ws4.Range("F20").Value = 23
$ws4.Range("F21").Value = 66
$ws4.Range("F25").Value = 26

$wb.Save()
